$d = $word.ActiveDocument

function Insert-ParaAfter([int]$index) {
    # Inserts a brand-new (empty) paragraph immediately after $d.Paragraphs($index)
    # using a range collapsed at that paragraph's end and InsertAfter("`r") --
    # this produces a clean <w:p/> (no stray empty run) with this engine,
    # unlike Range.InsertParagraphAfter()/InsertParagraphBefore().
    $p = $d.Paragraphs($index)
    $r = $d.Range($p.Range.End, $p.Range.End)
    $r.InsertAfter("`r")
}

# ---------------------------------------------------------------------------
# 1) Insert the new opening line "What has DY done for me lately? A lot"
#    followed by a blank paragraph, before the existing first paragraph.
# ---------------------------------------------------------------------------
$d.Paragraphs(1).Range.InsertParagraphBefore()
$d.Paragraphs(1).Range.Text = 'What has DY done for me lately? A lot'
Insert-ParaAfter 1

# ---------------------------------------------------------------------------
# 2) Extend the "Thank you all..." paragraph with the new closing clause.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    'they mean.', $true, $false, $false, $false, $false, $true, 1, $false,
    'they mean and zip through them rather expeditiously.', 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Rework the "This talk aims to briefly touch upon..." paragraph.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    'aims to briefly touch', $true, $false, $false, $false, $false, $true, 1, $false,
    'will touch', 2) | Out-Null

$d.Content.Find.Execute(
    'database. It’s larger aim', $true, $false, $false, $false, $false, $true, 1, $false,
    'database. Its larger aim', 2) | Out-Null

$d.Content.Find.Execute(
    'data model and methodology', $true, $false, $false, $false, $false, $true, 1, $false,
    'data model, and how our methodology', 2) | Out-Null

$d.Content.Find.Execute(
    'two major issues longevity', $true, $false, $false, $false, $false, $true, 1, $false,
    'two major issues: longevity', 2) | Out-Null

$d.Content.Find.Execute(
    'interoperability: projects tend to have a short-life cycle and DH projects are rarely able to speak to each other without significant data wrangling.',
    $true, $false, $false, $false, $false, $true, 1, $false,
    'interoperability.', 2) | Out-Null

$d.Content.Find.Execute(
    'For fear of sounding like I’m trying get you to buy a time-share, I will point out that  ',
    $true, $false, $false, $false, $false, $true, 1, $false,
    '', 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Insert the large block of brand-new paragraphs right after that
#    paragraph (which now ends "...longevity and interoperability. ").
#    (index-based; re-fetch $d.Paragraphs($idx) each time since stale
#    Paragraph object refs / .Next() can go stale once new paragraphs are
#    spliced into the document)
# ---------------------------------------------------------------------------
$idx = 4   # "This talk will touch upon..." paragraph, 1-based

$newParas = @(
    'DH projects have fractional life-cycles compared to print scholarship. Projects that are launched with great fan-fare often become defunct only a few years later. Solving the longevity question is not merely an academic question it is also a moral one. Let’s face it, DH programs suck up an incredible amount of funding. For 2023, the NEH announced it will disburse 2.2 million dollars in grants. For that amount of money, you could hire 30 “traditional” English professors who could churn out articles and books in print that won’t simply evaporate with the latest browser update. ',
    'Additionally, DH projects are rarely able to speak to each other without significant data wrangling. This is important because DH projects therefore do not benefit from the network effects visible in other disciplines such as political science, economics, or the hard sciences. There the data from individual projects can be recycled to other projects. For example, voting data from the US can be compared to voting data from Canada, even if there are important nuances in how that data is collected. Meanwhile, in the humanities, the data from a project about author X can rarely be used effectively as a point of comparison with author Y. Obviously, there are important exceptions to this such as the Seshat project, Stanford Standard Corpus, and others, but by and large this is the landscape. ',
    'In this sense, longevity and interoperability are intimately related. By making the data interoperable we are ensuring its reincarnation and upcycling in other projects. ',
    'I would argue that the data model created, refined, and tested for nearly a decade by the Digital Yoknapatawpha team represents a good starting point for a more universal framework. ',
    'Without going into too much technical detail, we can say that DY collects three main data enteties through close-reading texts: characters, locations, events. ',
    'The relationship between these three entities represents one of the fundamental structures of fiction: people (characters) doing something (events) somewhere (location). This flexible structure could arguably be applied and modified for most any author. Doing so would allow us to ask and answer questions that exceed a single author, and start to map the larger literary ecosystem. This in turn could help us reimagine the perennial question of this conference: What do we mean when we say Modernism?',
    'Lighting Example: Faulkner and Plot Structure',
    'One of the distinguish features of Faulkner, and, perhaps Modernism, is narrative experimentation. In our data, we have information that will let use visualize Faulkner’s use of time and narrative voice. '
)

foreach ($txt in $newParas) {
    Insert-ParaAfter $idx
    $idx = $idx + 1
    $d.Paragraphs($idx).Range.Text = $txt
}

# ---------------------------------------------------------------------------
# 5) Add three more blank paragraphs at the end (2 -> 5 total).
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt 3; $i++) {
    $n = $d.Paragraphs.Count
    Insert-ParaAfter $n
}

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
